$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (MT, 2015) - model changes from Exp to Sph, values updated
$ws.Range("C2").Value = "Sph"
$ws.Range("D2").Value = 0.111
$ws.Range("E2").Value = 1.564
$ws.Range("F2").Value = 4.47
$ws.Range("G2").Value = 0.07097186700767263
$ws.Range("H2").Value = 2.7726
$ws.Range("I2").Value = 0.7851

# Add new row 3 (MT, 2017)
$ws.Range("A3").Value = "MT"
$ws.Range("B3").Value = 2017
$ws.Range("C3").Value = "Exp"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2.0182
$ws.Range("F3").Value = 2.17
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 4.9353
$ws.Range("I3").Value = 0.0134

# Add new row 4 (MT, 2018)
$ws.Range("A4").Value = "MT"
$ws.Range("B4").Value = 2018
$ws.Range("C4").Value = "Gau"
$ws.Range("D4").Value = 0.8058999999999999
$ws.Range("E4").Value = 1.1195
$ws.Range("F4").Value = 4.2
$ws.Range("G4").Value = 0.7198749441715051
$ws.Range("H4").Value = 147.892
$ws.Range("I4").Value = 0.1663
